$d = $word.ActiveDocument

# Helper: append a new paragraph at the very end of the document whose
# text is exactly $text (may be empty), leaving a "clean" <w:p/> in the
# empty case instead of a stray empty run. We do this by inserting a
# unique marker run after the last paragraph, then using Find/Replace to
# turn "^p<marker>" into "^p<text>" in one shot — the trailing part of a
# Find/Replace that lands exactly on the (empty) tail of the document
# serializes without a leftover empty <w:r/>.
function AppendCleanParagraph {
    param([string]$Text)

    $lastIndex = $d.Paragraphs.Count
    $lastPara = $d.Paragraphs.Item($lastIndex)
    $tail = $lastPara.Range
    $tail.Collapse(0)
    $tail.InsertAfter("ZZ_NEW_PARA_MARKER_ZZ")

    $find = $d.Content.Find
    $find.Execute("ZZ_NEW_PARA_MARKER_ZZ", $false, $false, $false, $false, $false, `
                  $true, 1, $false, "^p" + $Text, 2)
}

# Remember where the new block starts so we can format the divider
# paragraph afterwards (formatting applied mid-stream would otherwise be
# inherited by every paragraph appended later).
$insertStart = $d.Paragraphs.Count

AppendCleanParagraph ""
AppendCleanParagraph "Vamos entrar em um novo repositório de alguém que deseja e clonar esse projeto."
AppendCleanParagraph ""
AppendCleanParagraph "Para que os projetos de outra pessoa possam ser criados como repositório no GitHub você precisará fazer um ‘Fork’."
AppendCleanParagraph ""
AppendCleanParagraph ""

# Give the first new paragraph (the blank divider right after "Commit")
# the same bottom-border formatting used by the divider under the title.
$dividerPara = $d.Paragraphs.Item($insertStart + 1)
$borders = $dividerPara.Borders
$borders.DistanceFromBottom = 1
$bottomBorder = $borders.Item(-3)
$bottomBorder.LineStyle = 1
$bottomBorder.LineWidth = 3
$bottomBorder.Color = -16777216
